{"js": "// Find the paragraph describing the team-member rights / materialized view\n// (identified by its distinctive opening text) and replace its content with\n// the revised wording from the commit: split sentences with extra commas,\n// and insert the new passage about the access-rights / log-creation\n// limitation before the \"Materializovan\u00fd pohled vznikl...\" sentence.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst marker = \"Druh\u00e9mu \u010dlenovi t\u00fdmu jsme p\u0159id\u011blili\";\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(marker) !== -1) {\n    target = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error(\"Target paragraph not found\");\n}\n\nconst newText =\n  \"Druh\u00e9mu \u010dlenovi t\u00fdmu jsme p\u0159id\u011blili v\u0161echna pr\u00e1va na v\u0161echny tabulky a procedury, kter\u00e9 jsme vytvo\u0159ili, co\u017e mu umo\u017enilo vytvo\u0159it materializovan\u00fd pohled z tabulek prvn\u00edho \u010dlena t\u00fdmu po tom, co byli prvn\u00edm \u010dlenem vytvo\u0159en\u00e9 logy. Logy nem\u016f\u017ee vytvo\u0159it druh\u00fd \u010dlen, proto\u017ee ani jeden z n\u00e1s nem\u00e1 dostate\u010dn\u00e1 opr\u00e1vn\u011bn\u00ed na p\u0159id\u011blen\u00ed pr\u00e1v na vytvo\u0159en\u00ed log\u016f.  Materializovan\u00fd pohled vznikl spojen\u00edm tabulek Letenka a Let, co\u017e p\u0159edstavuje informace o\\u00A0jednotliv\u00fdch cestuj\u00edc\u00edch na jednotliv\u00fdch letech. N\u00e1sledn\u011b jsme p\u0159id\u011blili ve\u0161ker\u00e1 p\u0159\u00edstupov\u00e1 pr\u00e1va na dan\u00fd materializovan\u00fd pohled prvn\u00edmu \u010dlenu t\u00fdmu.\";\n\ntarget.insertText(newText, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Locate the paragraph describing the access rights / materialized view\n# by its distinctive opening text.\n$target = $null\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Text -like \"*Druh\u00e9mu \u010dlenovi t\u00fdmu jsme p\u0159id\u011blili*\") {\n        $target = $p\n        break\n    }\n}\n\nif ($target -eq $null) {\n    throw \"Target paragraph not found\"\n}\n\n$newText = \"Druh\u00e9mu \u010dlenovi t\u00fdmu jsme p\u0159id\u011blili v\u0161echna pr\u00e1va na v\u0161echny tabulky a procedury, kter\u00e9 jsme vytvo\u0159ili, co\u017e mu umo\u017enilo vytvo\u0159it materializovan\u00fd pohled z tabulek prvn\u00edho \u010dlena t\u00fdmu po tom, co byli prvn\u00edm \u010dlenem vytvo\u0159en\u00e9 logy. Logy nem\u016f\u017ee vytvo\u0159it druh\u00fd \u010dlen, proto\u017ee ani jeden z n\u00e1s nem\u00e1 dostate\u010dn\u00e1 opr\u00e1vn\u011bn\u00ed na p\u0159id\u011blen\u00ed pr\u00e1v na vytvo\u0159en\u00ed log\u016f.  Materializovan\u00fd pohled vznikl spojen\u00edm tabulek Letenka a Let, co\u017e p\u0159edstavuje informace o\u00a0jednotliv\u00fdch cestuj\u00edc\u00edch na jednotliv\u00fdch letech. N\u00e1sledn\u011b jsme p\u0159id\u011blili ve\u0161ker\u00e1 p\u0159\u00edstupov\u00e1 pr\u00e1va na dan\u00fd materializovan\u00fd pohled prvn\u00edmu \u010dlenu t\u00fdmu.\"\n\n# Keep the paragraph mark (and trailing bookmark) intact: shrink the range\n# by one character before overwriting its text.\n$r = $target.Range\n$r.MoveEnd(1, -1)\n$r.Text = $newText\n"}
